# Replaces the alert/part-number table (rows 2..52) with the corrected
# dataset, then re-applies the conditional formatting over the new extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: PART NO. index, part number, value, alert color
$data = @(
    @(0,  "CPA110Y-NMSNN-NNNNN", 0,     "RED"),
    @(1,  "D0114RB",             668,   "BLUE"),
    @(2,  "D0117XL-A",           950,   "BLUE"),
    @(3,  "F9270BJ",             1968,  "BLUE"),
    @(4,  "F9270BK",             2002,  "BLUE"),
    @(5,  "F9273CZ",             3798,  "BLUE"),
    @(6,  "F9300AF",             3027,  "BLUE"),
    @(7,  "F9300TE",             927,   "BLUE"),
    @(8,  "F9340XX",             656,   "BLUE"),
    @(9,  "F9341JA",             2031,  "BLUE"),
    @(10, "F9341JD",             2031,  "BLUE"),
    @(11, "F9341JE",             2031,  "BLUE"),
    @(12, "F9341JP",             1996,  "BLUE"),
    @(13, "F9341JQ",             2028,  "BLUE"),
    @(14, "F9341JS",             2002,  "BLUE"),
    @(15, "F9341JW",             1096,  "BLUE"),
    @(16, "F9341JZ",             1089,  "BLUE"),
    @(17, "F9342MK",             3975,  "BLUE"),
    @(18, "F9342NJ",             2000,  "BLUE"),
    @(19, "F9342YB",             475,   "BLUE"),
    @(20, "F9900AC",             2276,  "BLUE"),
    @(21, "F9900AU",             860,   "BLUE"),
    @(22, "F9900BC",             2276,  "BLUE"),
    @(23, "F9900EA",             2096,  "BLUE"),
    @(24, "F9900GC",             3936,  "BLUE"),
    @(25, "F9900GE",             1838,  "BLUE"),
    @(26, "F9900GH",             1996,  "BLUE"),
    @(27, "F9900GZ",             6795,  "BLUE"),
    @(28, "F9900RE",             2009,  "BLUE"),
    @(29, "F9900RF",             1128,  "BLUE"),
    @(30, "F9900RG",             10093, "BLUE"),
    @(31, "F9900RJ",             2002,  "BLUE"),
    @(32, "F9900RS",             2009,  "BLUE"),
    @(33, "F9900TB",             1522,  "BLUE"),
    @(34, "F9903BD",             85,    "BLUE"),
    @(35, "F9903CM",             83,    "BLUE"),
    @(36, "F9903CZ",             83,    "BLUE"),
    @(37, "F9910LE",             1996,  "BLUE"),
    @(38, "F9910LS",             2002,  "BLUE"),
    @(39, "F9913QG",             777,   "BLUE"),
    @(40, "F9913QH",             983,   "BLUE"),
    @(41, "F9913TA",             1805,  "BLUE"),
    @(42, "F9913TS",             1851,  "BLUE"),
    @(43, "F9920LU",             240,   "BLUE"),
    @(44, "F9921AH",             159,   "BLUE"),
    @(45, "F9921VG",             25,    "BLUE"),
    @(46, "G9303NC",             1996,  "BLUE"),
    @(47, "G9330DB",             4700,  "BLUE"),
    @(48, "Y9308JY",             4004,  "BLUE"),
    @(49, "Y9501WL",             6006,  "BLUE"),
    @(50, "Y9812PS",             2002,  "BLUE")
)

$oldLastRow = 35
$newLastRow = 1 + $data.Count   # header is row 1

# Clear out any previous data rows beyond the new extent (none expected here,
# but keep the sheet tidy if the old table was ever longer).
if ($oldLastRow -gt $newLastRow) {
    $ws.Range("A$($newLastRow + 1):D$oldLastRow").Clear() | Out-Null
}

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $entry = $data[$i]
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
}

# Column A carries the bold/bordered "s=1" style on every data row. Rows
# beyond the old table extent are brand-new cells, so copy that formatting
# down, one source cell at a time, from the existing column-A cells (values
# are untouched - formats only; sizes must match exactly so the paste isn't
# tiled past the intended destination row).
if ($newLastRow -gt $oldLastRow) {
    for ($row = $oldLastRow + 1; $row -le $newLastRow; $row++) {
        $ws.Range("A2").Copy() | Out-Null
        $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
    }
    $excel.CutCopyMode = 0
}

# Re-point the conditional formatting at the full, now-larger table range
# (keeps the existing rules / dxf styles intact, only the sqref changes).
$newRange = $ws.Range("A1:D$newLastRow")
foreach ($fc in $ws.Range("A1:D$oldLastRow").FormatConditions) {
    $fc.ModifyAppliesToRange($newRange)
}
